# "file update to 28 april" - append the 28-Apr-2020 (serial 43949) daily
# figures as a new row (row 53) at the bottom of each of the three
# worksheets (Confirmed, Recoverd, Death). Column B keeps the running
# "cumulative total" formula (=SUM(previous B + this row's C)), column C
# holds the day's new count.

$wb = $excel.ActiveWorkbook

# Per-sheet data for 28-Apr-2020 and the new selected cell left behind
# on that sheet after the edit.
$sheetInfo = @(
    @{ Name = "Confirmed"; C = 549; Selection = "D55" },
    @{ Name = "Recoverd";  C = 8;   Selection = "B53" },
    @{ Name = "Death";     C = 3;   Selection = "G38" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # New date cell, matching the style/format used by the rest of column A.
    $ws.Range("A53").Value = 43949
    $ws.Range("A53").NumberFormat = "d-mmm-yy"

    # New "new cases" cell, matching column C's center/center style.
    $ws.Range("C53").Value = $info.C

    # Running total formula, continuing the pattern already used down
    # column B (=SUM(B<prevRow>+C<thisRow>)).
    $ws.Range("B53").Formula = "=SUM(B52+C53)"

    # Match the center horizontal/vertical alignment style (s="2") used
    # by every other data cell in columns B and C.
    $ws.Range("B53:C53").HorizontalAlignment = -4108
    $ws.Range("B53:C53").VerticalAlignment = -4108

    # Move the selection to wherever the author left it after entering
    # the new row (the view's frozen/scrolled "topLeftCell" is restored
    # automatically by the host application when the file is reopened).
    $ws.Range($info.Selection).Select()
}
